# Update "想去人数" (want-to-go count) figures to the values captured at
# generation time 456a3b4 (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 26917
$ws.Range("F4").Value = 601
$ws.Range("F8").Value = 558
$ws.Range("F11").Value = 457
$ws.Range("F19").Value = 227
$ws.Range("F20").Value = 64
$ws.Range("F22").Value = 106

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5126

# 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5126
$ws.Range("F5").Value = 26917
$ws.Range("F6").Value = 601
$ws.Range("F19").Value = 558
$ws.Range("F23").Value = 457
$ws.Range("F35").Value = 227
$ws.Range("F37").Value = 64
$ws.Range("F39").Value = 106
